# Trade #75 closed at 2026-02-17 08:58:06 - unknown UNKNOWN +0.000%
#
# Updates the "Summary" and "Strategy Status" aggregate figures for the
# MarketMaking strategy, and appends the newly-closed trade row (#75) to
# both the "All Trades" and "MarketMaking" trade logs.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Summary sheet - top level account metrics
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.52   # Current Capital
$summary.Range("B4").Value = 0.53      # Total P&L $
$summary.Range("B5").Value = 0.14      # Total P&L %
$summary.Range("B6").Value = 75        # Total Trades
$summary.Range("B8").Value = 30        # Losing Trades
$summary.Range("B9").Value = 41.33     # Win Rate %

# ---------------------------------------------------------------------
# 2. Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.52     # Capital
$status.Range("D4").Value = 75         # Trades
$status.Range("E4").Value = 0.53       # P&L $
$status.Range("F4").Value = 0.52       # P&L %
$status.Range("G4").Value = 41.33      # Win Rate %

# ---------------------------------------------------------------------
# 3. Append the new closed trade (#75) to the trade logs
# ---------------------------------------------------------------------
function Add-TradeRow($ws, [int]$row) {
    $ws.Range("A$row").Value = 75
    # Leading apostrophe forces text entry so Excel does not reinterpret
    # the date-shaped string as a date serial number (matches the source
    # file, which stores these as plain inline strings). Reset the style
    # afterwards so the cell keeps the workbook's default formatting
    # instead of picking up a "quote prefix" indicator.
    $ws.Range("B$row").Value = "'2026-02-17"
    $ws.Range("B$row").Style = "Normal"
    $ws.Range("C$row").Value = "08:58:00"
    $ws.Range("D$row").Value = "MarketMaking"
    $ws.Range("E$row").Value = "DOWN"
    $ws.Range("F$row").Value = 0.91
    $ws.Range("G$row").Value = 0.86
    $ws.Range("H$row").Value = "CLOSED"
    $ws.Range("I$row").Value = -5.4945
    $ws.Range("J$row").Value = -0.05
    $ws.Range("K$row").Value = 100.52
    $ws.Range("L$row").Value = 0
    $ws.Range("M$row").Value = 0
    $ws.Range("N$row").Value = 0.6
    $ws.Range("O$row").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P$row").Value = "early_exit"
    $ws.Range("Q$row").Value = 0.13
}

$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $allTrades 76

$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $marketMaking 76
